# Turn the three contact-info values (email, GitHub, LinkedIn) into real
# hyperlinks, matching how Word's AutoFormat / Insert Hyperlink produces
# a <w:hyperlink> wrapping a run styled with the built-in "Hyperlink"
# character style.

$d = $word.ActiveDocument

# Email address -> mailto: link
$rEmail = $d.Content.Duplicate
$rEmail.Find.Execute("JosephTLyons@gmail.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($rEmail, "mailto:JosephTLyons@gmail.com", "", "", "JosephTLyons@gmail.com")

# GitHub profile URL (the "www." + "github.com/..." runs collapse into the
# single run that now lives inside the hyperlink)
$rGitHub = $d.Content.Duplicate
$rGitHub.Find.Execute("www.github.com/JosephTLyons", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($rGitHub, "http://www.github.com/JosephTLyons", "", "", "www.github.com/JosephTLyons")

# LinkedIn profile URL
$rLinkedIn = $d.Content.Duplicate
$rLinkedIn.Find.Execute("www.linkedin.com/in/JosephTLyons", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($rLinkedIn, "http://www.linkedin.com/in/JosephTLyons", "", "", "www.linkedin.com/in/JosephTLyons")
